$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos table
# with the latest scraped quotes. Values are written with a leading
# apostrophe (quote-prefix) so Excel keeps them as literal text (matching
# the source data, e.g. "27.317.73", "0.07700", "  -0.64%  ") instead of
# silently re-interpreting/normalizing them as numbers.
$ws.Range('D2').Value = '''27.317.73'
$ws.Range('E2').Value = '''  -0.64%  '
$ws.Range('D3').Value = '''1.716.73'
$ws.Range('E4').Value = '''  -0.09%  '
$ws.Range('D5').Value = '''225.55'
$ws.Range('E5').Value = '''  +0.07%  '
$ws.Range('D6').Value = '''0.5306'
$ws.Range('D7').Value = '''1.003'
$ws.Range('E7').Value = '''  -0.09%  '
$ws.Range('D8').Value = '''0.06703'
$ws.Range('E8').Value = '''  +1.82%  '
$ws.Range('D9').Value = '''0.2668'
$ws.Range('E9').Value = '''  +0.15%  '
$ws.Range('D10').Value = '''20.96'
$ws.Range('E10').Value = '''  -2.88%  '
$ws.Range('D11').Value = '''0.07700'
$ws.Range('E11').Value = '''  +0.04%  '
$ws.Range('D12').Value = '''4.511'
$ws.Range('E12').Value = '''  -1.96%  '
$ws.Range('D13').Value = '''1.952.67'
$ws.Range('E13').Value = '''  -0.51%  '
$ws.Range('D14').Value = '''1.725.46'
$ws.Range('E14').Value = '''  +0.02%  '
$ws.Range('D15').Value = '''0.5870'
$ws.Range('E15').Value = '''  +0.96%  '
$ws.Range('D16').Value = '''0.0₅8224'
$ws.Range('E16').Value = '''  -0.55%  '
$ws.Range('D17').Value = '''68.10'
$ws.Range('E17').Value = '''  +0.45%  '
$ws.Range('D18').Value = '''27.415.43'
$ws.Range('E18').Value = '''  -0.32%  '
$ws.Range('D19').Value = '''223.38'
$ws.Range('E19').Value = '''  +2.52%  '
$ws.Range('D20').Value = '''1.003'
$ws.Range('E20').Value = '''  -0.14%  '
$ws.Range('D21').Value = '''4.664'
$ws.Range('E21').Value = '''  -1.17%  '
$ws.Range('D23').Value = '''6.048'
$ws.Range('E23').Value = '''  -0.44%  '
$ws.Range('D24').Value = '''1.004'
$ws.Range('E24').Value = '''  -0.03%  '
$ws.Range('D25').Value = '''144.37'
$ws.Range('E25').Value = '''  +0.76%  '
$ws.Range('D26').Value = '''1.695'
$ws.Range('E26').Value = '''  -3.26%  '
$ws.Range('E27').Value = '''  -1.82%  '
$ws.Range('D28').Value = '''7.260'
$ws.Range('E28').Value = '''  -1.77%  '
$ws.Range('D29').Value = '''16.30'
$ws.Range('E29').Value = '''  -1.28%  '
$ws.Range('D30').Value = '''0.05369'
$ws.Range('E30').Value = '''  -2.11%  '
$ws.Range('D31').Value = '''1.296'
$ws.Range('E31').Value = '''  -0.48%  '
$ws.Range('D32').Value = '''3.488'
$ws.Range('E32').Value = '''  -2.07%  '
$ws.Range('D33').Value = '''3.435'
$ws.Range('E33').Value = '''  -0.14%  '
$ws.Range('D34').Value = '''1.631'
$ws.Range('E34').Value = '''  -1.61%  '
$ws.Range('D35').Value = '''2.870'
$ws.Range('E35').Value = '''  +0.28%  '
$ws.Range('D36').Value = '''0.9587'
$ws.Range('E36').Value = '''  -0.47%  '
$ws.Range('D37').Value = '''2.394'
$ws.Range('E37').Value = '''  -1.22%  '
$ws.Range('D38').Value = '''0.5883'
$ws.Range('E38').Value = '''  -1.39%  '
$ws.Range('D39').Value = '''1.148.66'
$ws.Range('E39').Value = '''  +9.07%  '
$ws.Range('D40').Value = '''0.01649'
$ws.Range('E40').Value = '''  +0.02%  '
$ws.Range('D41').Value = '''5.806'
$ws.Range('E41').Value = '''  -1.66%  '
$ws.Range('D42').Value = '''1.004'
$ws.Range('E42').Value = '''  +0.02%  '
$ws.Range('D43').Value = '''0.8431'
$ws.Range('E43').Value = '''  -0.92%  '
$ws.Range('D44').Value = '''100.97'
$ws.Range('E44').Value = '''  -0.36%  '
$ws.Range('D45').Value = '''1.858.64'
$ws.Range('E45').Value = '''  -0.59%  '
$ws.Range('E46').Value = '''  -5.54%  '
$ws.Range('D47').Value = '''57.79'
$ws.Range('E47').Value = '''  -1.75%  '
$ws.Range('D48').Value = '''0.4593'
$ws.Range('E48').Value = '''  +2.64%  '
$ws.Range('E49').Value = '''  +0.14%  '
$ws.Range('D50').Value = '''8.134'
$ws.Range('E50').Value = '''  -0.48%  '
$ws.Range('D51').Value = '''0.05197'
$ws.Range('E51').Value = '''  -0.79%  '
